$wb = $excel.ActiveWorkbook

# --- Update "Resultados" sheet (sheet1): refreshed GB column (D) and G-mean column (M) values ---
$ws1 = $wb.Worksheets.Item("Resultados")
$ws1.Range("D3").Value = 0.8445743889623244
$ws1.Range("M3").Value = 0.8238599096663985
$ws1.Range("D4").Value = 0.9442628728875527
$ws1.Range("M4").Value = 0.9185091437335826
$ws1.Range("D5").Value = 0.9442628728875527
$ws1.Range("M5").Value = 0.9217255693426315

# --- Translate "Selected Features" sheet (sheet3) from Portuguese to English ---
$ws3 = $wb.Worksheets.Item("Selected Features")
$ws3.Range('B2').Value = 'Neurological disorder/symptoms '
$ws3.Range('B3').Value = 'Neurological disorder/symptoms '
$ws3.Range('C3').Value = 'Difficulty walking'
$ws3.Range('B4').Value = 'Neurological disorder/symptoms '
$ws3.Range('C4').Value = 'Neuropathy'
$ws3.Range('D4').Value = 'Difficulty walking'
$ws3.Range('B5').Value = 'Sex'
$ws3.Range('C5').Value = 'Neurological disorder/symptoms '
$ws3.Range('D5').Value = 'Neuropathy'
$ws3.Range('E5').Value = 'Difficulty walking'
$ws3.Range('B6').Value = 'Sex'
$ws3.Range('C6').Value = 'Urinary symptoms'
$ws3.Range('D6').Value = 'Neurological disorder/symptoms '
$ws3.Range('E6').Value = 'Neuropathy'
$ws3.Range('F6').Value = 'Difficulty walking'
$ws3.Range('B7').Value = 'Sex'
$ws3.Range('C7').Value = 'Urinary symptoms'
$ws3.Range('D7').Value = 'Pain'
$ws3.Range('E7').Value = 'Neurological disorder/symptoms '
$ws3.Range('F7').Value = 'Neuropathy'
$ws3.Range('G7').Value = 'Difficulty walking'
$ws3.Range('B8').Value = 'Sex'
$ws3.Range('C8').Value = 'Age'
$ws3.Range('D8').Value = 'Urinary symptoms'
$ws3.Range('E8').Value = 'Pain'
$ws3.Range('F8').Value = 'Neurological disorder/symptoms '
$ws3.Range('G8').Value = 'Neuropathy'
$ws3.Range('H8').Value = 'Difficulty walking'
$ws3.Range('B9').Value = 'Sex'
$ws3.Range('C9').Value = 'Age'
$ws3.Range('D9').Value = 'Urinary symptoms'
$ws3.Range('E9').Value = 'Pain'
$ws3.Range('F9').Value = 'Neurological disorder/symptoms '
$ws3.Range('G9').Value = 'Neuropathy'
$ws3.Range('H9').Value = 'Difficulty walking'
$ws3.Range('I9').Value = 'Sphincter dyscontrol'
$ws3.Range('B10').Value = 'Sex'
$ws3.Range('C10').Value = 'Age'
$ws3.Range('E10').Value = 'Urinary symptoms'
$ws3.Range('F10').Value = 'Pain'
$ws3.Range('G10').Value = 'Neurological disorder/symptoms '
$ws3.Range('H10').Value = 'Neuropathy'
$ws3.Range('I10').Value = 'Difficulty walking'
$ws3.Range('J10').Value = 'Sphincter dyscontrol'
$ws3.Range('B11').Value = 'Sex'
$ws3.Range('C11').Value = 'Age'
$ws3.Range('E11').Value = 'Urinary symptoms'
$ws3.Range('F11').Value = 'Pain'
$ws3.Range('G11').Value = 'Neurological disorder/symptoms '
$ws3.Range('H11').Value = 'Neuropathy'
$ws3.Range('I11').Value = 'Mental disorders'
$ws3.Range('J11').Value = 'Difficulty walking'
$ws3.Range('K11').Value = 'Sphincter dyscontrol'
$ws3.Range('B12').Value = 'Sex'
$ws3.Range('C12').Value = 'Age'
$ws3.Range('E12').Value = 'Urinary symptoms'
$ws3.Range('F12').Value = 'Pain'
$ws3.Range('G12').Value = 'Neurological disorder/symptoms '
$ws3.Range('H12').Value = 'Neuropathy'
$ws3.Range('I12').Value = 'Chronic disease'
$ws3.Range('J12').Value = 'Mental disorders'
$ws3.Range('K12').Value = 'Difficulty walking'
$ws3.Range('L12').Value = 'Sphincter dyscontrol'
$ws3.Range('B13').Value = 'Sex'
$ws3.Range('C13').Value = 'Age'
$ws3.Range('F13').Value = 'Urinary symptoms'
$ws3.Range('G13').Value = 'Pain'
$ws3.Range('H13').Value = 'Neurological disorder/symptoms '
$ws3.Range('I13').Value = 'Neuropathy'
$ws3.Range('J13').Value = 'Chronic disease'
$ws3.Range('K13').Value = 'Mental disorders'
$ws3.Range('L13').Value = 'Difficulty walking'
$ws3.Range('M13').Value = 'Sphincter dyscontrol'
$ws3.Range('B14').Value = 'Sex'
$ws3.Range('C14').Value = 'Age'
$ws3.Range('G14').Value = 'Urinary symptoms'
$ws3.Range('H14').Value = 'Pain'
$ws3.Range('I14').Value = 'Neurological disorder/symptoms '
$ws3.Range('J14').Value = 'Neuropathy'
$ws3.Range('K14').Value = 'Chronic disease'
$ws3.Range('L14').Value = 'Mental disorders'
$ws3.Range('M14').Value = 'Difficulty walking'
$ws3.Range('N14').Value = 'Sphincter dyscontrol'
